$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.965.99'
$ws.Range("E2").Value = '  -6.14%  '

$ws.Range("D3").Value = '2.884.09'
$ws.Range("E3").Value = '  -3.67%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '546.00'
$ws.Range("E5").Value = '  -2.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '123.22'

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.499'
$ws.Range("E8").Value = '  +0.92%  '

$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '2.885.04'
$ws.Range("E9").Value = '  -3.56%  '

$ws.Range("E10").Value = '  -9.63%  '

$ws.Range("E11").Value = '  +0.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.62'
$ws.Range("E12").Value = '  -10.62%  '

$ws.Range("E13").Value = '  -5.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.46'
$ws.Range("E14").Value = '  -1.13%  '

$ws.Range("E15").Value = '  +1.10%  '

$ws.Range("D16").Value = '3.360.43'
$ws.Range("E16").Value = '  -3.61%  '

$ws.Range("D17").Value = '2.884.16'
$ws.Range("E17").Value = '  -3.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.51'
$ws.Range("E18").Value = '  +5.85%  '

$ws.Range("D19").Value = '57.019.81'
$ws.Range("E19").Value = '  -6.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '404.02'
$ws.Range("E20").Value = '  -6.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.78'
$ws.Range("E21").Value = '  -2.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.670'
$ws.Range("E22").Value = '  +1.55%  '

$ws.Range("E23").Value = '  -4.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.66'
$ws.Range("E24").Value = '  -1.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '77.26'
$ws.Range("E25").Value = '  -1.85%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("E28").Value = '  -2.00%  '

$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.92'
$ws.Range("E29").Value = '  +2.59%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.17'
$ws.Range("E30").Value = '  +0.25%  '

$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("E32").Value = '  -3.06%  '

$ws.Range("E33").Value = '  +5.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.905'
$ws.Range("E34").Value = '  -5.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.37'
$ws.Range("E35").Value = '  -3.58%  '

$ws.Range("E36").Value = '  -11.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '47.75'
$ws.Range("E37").Value = '  -4.25%  '

$ws.Range("E39").Value = '  -6.91%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.105'
$ws.Range("E40").Value = '  -1.60%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0337'
$ws.Range("E41").Value = '  -5.78%  '

$ws.Range("D42").Value = '2.625.21'
$ws.Range("E42").Value = '  -1.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '357.03'
$ws.Range("E43").Value = '  -4.40%  '

$ws.Range("E44").Value = '  -2.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '119.60'
$ws.Range("E46").Value = '  -0.71%  '

$ws.Range("E47").Value = '  -2.71%  '

$ws.Range("E48").Value = '  +0.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.92'
$ws.Range("E49").Value = '  -1.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.70'
$ws.Range("E50").Value = '  -2.61%  '

$ws.Range("E51").Value = '  -3.85%  '
